# "Add files via upload" — append two new "word list" rows at the bottom
# of Sheet1: A91 = "por" (an existing word, reusing the shared string),
# A92 = "paí" (a brand-new word), styled with the Roboto / #252525 font
# used for this new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dimension / used range grows from A1:A90 to A1:A92
# automatically once these are written).
$ws.Range("A91").Value = "por"
$ws.Range("A92").Value = "paí"

# A92 gets its own font: Roboto, 11pt, RGB(0x25,0x25,0x25).
$cell = $ws.Range("A92")
$cell.Font.Name = "Roboto"
$cell.Font.Size = 11
$cell.Font.Color = 2434341   # RGB(37,37,37) = 0x252525

# Reflect the new bottom-of-sheet position in the sheet view.
$cell.Select() | Out-Null

# Printable page setup (A4 portrait), matching the refreshed file.
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
